{"js": "// Add \"The proposed solution can meet project timeline\" as a new bullet\n// right before the existing \"Architecture design can meet the business use\n// cases or end user requirements\" bullet, and tighten the wording of the\n// Kibana bullet (\"lot of\" -> \"additional\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Architecture design ...\" bullet \u2014 this is the paragraph whose\n// text the new \"project timeline\" bullet must be inserted immediately before.\nlet architectureParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Architecture design can meet the business use cases\") !== -1) {\n    architectureParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!architectureParagraph) {\n  throw new Error('Could not find the \"Architecture design can meet the business use cases\" paragraph.');\n}\n\n// Insert the new bullet before it; insertParagraph inherits the paragraph\n// (pPr) and run (rPr) formatting of the reference paragraph automatically.\narchitectureParagraph.insertParagraph(\n  \"The proposed solution can meet project timeline\",\n  \"Before\"\n);\nawait context.sync();\n\n// Locate the Kibana bullet and tighten \"lot of BI tools\" -> \"additional BI tools\".\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet kibanaParagraph = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"Kibana is the only BI tool\") !== -1) {\n    kibanaParagraph = paragraphs2.items[i];\n    break;\n  }\n}\n\nif (!kibanaParagraph) {\n  throw new Error(\"Could not find the Kibana paragraph.\");\n}\n\nconst kibanaRange = kibanaParagraph.getRange();\nconst hits = kibanaRange.search(\"lot of\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items/text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"lot of\" inside the Kibana paragraph.');\n}\n\nhits.items[0].insertText(\"additional\", \"Replace\");\nawait context.sync();\n", "ps1": "# Add \"The proposed solution can meet project timeline\" as a new bullet\n# right before the existing \"Architecture design can meet the business use\n# cases or end user requirements\" bullet, and tighten the wording of the\n# Kibana bullet (\"lot of\" -> \"additional\").\n\n$d = $word.ActiveDocument\n\n# Locate the \"Architecture design ...\" bullet paragraph.\n$architectureParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Architecture design can meet the business use cases*\") {\n        $architectureParagraph = $p\n        break\n    }\n}\n\nif ($null -eq $architectureParagraph) {\n    throw 'Could not find the \"Architecture design can meet the business use cases\" paragraph.'\n}\n\n# Insert a new paragraph right before it (inherits its pPr/rPr formatting),\n# then fill it in with the new assumption text.\n$insertionPoint = $architectureParagraph.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n$newParagraphRange = $insertionPoint.InsertParagraphBefore()\n\n$d = $word.ActiveDocument\n$newParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Architecture design can meet the business use cases*\") {\n        $newParagraph = $p.Previous()\n        break\n    }\n}\n$newParagraph.Range.Text = \"The proposed solution can meet project timeline\"\n\n# Locate the Kibana bullet and tighten \"lot of BI tools\" -> \"additional BI tools\".\n$kibanaParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Kibana is the only BI tool*\") {\n        $kibanaParagraph = $p\n        break\n    }\n}\n\nif ($null -eq $kibanaParagraph) {\n    throw \"Could not find the Kibana paragraph.\"\n}\n\n$find = $kibanaParagraph.Range.Find\n$null = $find.Execute(\"lot of\", $false, $false, $false, $false, $false, $true, 1, $false, \"additional\", 2)\n"}
